$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6, pushing existing rows 6-34 down to 7-35.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record.
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44462
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 380
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
